$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 397; this pushes the existing rows 397..439
# down to 398..440 (preserving all their data/formatting).
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with the new record.
$ws.Range("A397").Value = 4
$ws.Range("B397").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C397").Value = 'Los Lagos'
$ws.Range("D397").Value = 45124
$ws.Range("E397").Value = 10
$ws.Range("F397").Value = 'Fruta'
$ws.Range("G397").Value = 100108
$ws.Range("H397").Value = 'Tropicales y subtropicales'
$ws.Range("I397").Value = 100108005
$ws.Range("J397").Value = 'Piña'
$ws.Range("K397").Value = 'Caramelo'
$ws.Range("L397").Value = 'Primera'
$ws.Range("M397").Value = 30
$ws.Range("N397").Value = 27000
$ws.Range("O397").Value = 27000
$ws.Range("P397").Value = 27000
$ws.Range("Q397").Value = '$/caja 12 unidades'
$ws.Range("R397").Value = 'Ecuador'
$ws.Range("S397").Value = 2250
$ws.Range("T397").Value = 12
